$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-06 14:14:37"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-06 14:14:12"
$wsZhCn.Range("K4").Value = "2016-09-06 14:15:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-06 14:14:37"
$wsDeDe.Range("K4").Value = "2016-09-06 14:16:35"
